$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with revised M2 figures (open/high/low/close) ---

# Row 434 (2021-01-01): 1183618000000 -> 1182578000000
$ws.Range("C434:F434").Value = 1182578000000

# Row 435 (2021-02-01): 1181166000000 -> 1180145000000
$ws.Range("C435:F435").Value = 1180145000000

# Row 451 (2022-06-01): 1297397000000 -> 1297962000000
$ws.Range("C451:F451").Value = 1297962000000

# Row 457 (2022-12-01): 1376065000000 -> 1375618000000
$ws.Range("C457:F457").Value = 1375618000000

# Row 458 (2023-01-01): 1365331000000 -> 1365405000000
$ws.Range("C458:F458").Value = 1365405000000

# --- Append new rows 460-462 with the same look & feel as row 459 ---

$ws.Range("A459:G459").Copy($ws.Range("A460:G460"))
$ws.Cells.Item(460, 1).Value = 44986.45833333334
$ws.Cells.Item(460, 2).Value = "ECONOMICS:MAM2"
$ws.Range("C460:F460").Value = 1393379000000
$ws.Cells.Item(460, 7).Value = 0

$ws.Range("A459:G459").Copy($ws.Range("A461:G461"))
$ws.Cells.Item(461, 1).Value = 45017.45833333334
$ws.Cells.Item(461, 2).Value = "ECONOMICS:MAM2"
$ws.Range("C461:F461").Value = 1391394000000
$ws.Cells.Item(461, 7).Value = 0

$ws.Range("A459:G459").Copy($ws.Range("A462:G462"))
$ws.Cells.Item(462, 1).Value = 45047.41666666666
$ws.Cells.Item(462, 2).Value = "ECONOMICS:MAM2"
$ws.Range("C462:F462").Value = 1394414000000
$ws.Cells.Item(462, 7).Value = 0
